$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the numeric strings ("20000","10000","5000") from the shared-string pool by
# overwriting column B with real numbers, then add a selection to the sheet view.

# Block A=5 (B=20000): rows 2-32, reversed order
$rows5 = @(@(28,1), @(82,2), @(52,3), @(77,4), @(50,6), @(26,13), @(57,17), @(5,19), @(37,22), @(14,23), @(3,25), @(61,28), @(70,29), @(67,30), @(89,31), @(60,32), @(79,33), @(15,34), @(51,36), @(68,42), @(47,52), @(49,57), @(33,59), @(34,69), @(58,70), @(38,74), @(54,75), @(56,76), @(12,79), @(8,82), @(22,85))
for ($i = 0; $i -lt $rows5.Length; $i++) {
    $r = 2 + $i
    $pair = $rows5[$i]
    $ws.Cells.Item($r, 1).Value = 5
    $ws.Cells.Item($r, 2).Value = 20000
    $ws.Cells.Item($r, 3).Value = $pair[0]
    $ws.Cells.Item($r, 4).Value = $pair[1]
}

# Block A=2 (B=10000): rows 33-76, reversed order
$rows2 = @(@(28,1), @(82,2), @(52,3), @(1,7), @(25,11), @(36,16), @(5,19), @(31,20), @(14,23), @(60,32), @(59,35), @(73,38), @(20,39), @(30,43), @(69,44), @(65,45), @(75,46), @(83,47), @(87,49), @(47,52), @(81,58), @(94,60), @(19,61), @(16,62), @(62,63), @(11,64), @(35,65), @(2,66), @(76,68), @(72,71), @(10,72), @(54,75), @(56,76), @(71,77), @(21,78), @(24,80), @(74,81), @(8,82), @(32,83), @(29,84), @(7,86), @(64,87), @(53,88), @(44,89))
for ($i = 0; $i -lt $rows2.Length; $i++) {
    $r = 33 + $i
    $pair = $rows2[$i]
    $ws.Cells.Item($r, 1).Value = 2
    $ws.Cells.Item($r, 2).Value = 10000
    $ws.Cells.Item($r, 3).Value = $pair[0]
    $ws.Cells.Item($r, 4).Value = $pair[1]
}

# Block A=1 (B=5000): rows 77-121, reversed order
$rows1 = @(@(28,1), @(82,2), @(77,4), @(55,5), @(50,6), @(17,8), @(48,9), @(46,10), @(25,11), @(63,12), @(26,13), @(4,14), @(41,15), @(57,17), @(13,18), @(45,21), @(37,22), @(14,23), @(23,24), @(6,26), @(86,27), @(61,28), @(88,37), @(73,38), @(84,40), @(40,41), @(69,44), @(83,47), @(78,48), @(87,49), @(9,50), @(18,51), @(85,53), @(39,54), @(43,55), @(80,56), @(49,57), @(81,58), @(33,59), @(19,61), @(62,63), @(2,66), @(42,67), @(27,73), @(54,75))
for ($i = 0; $i -lt $rows1.Length; $i++) {
    $r = 77 + $i
    $pair = $rows1[$i]
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = 5000
    $ws.Cells.Item($r, 3).Value = $pair[0]
    $ws.Cells.Item($r, 4).Value = $pair[1]
}

$ws.Range("A1:D121").Select()
